# Append the new Webstaurant Bakery order lines (rows 9-17) to the sheet.
# All existing data cells (including the numeric-looking quantity/cost
# columns) are stored as text, so force the Text number format before
# writing any value to keep that convention for the new rows too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("612PIE8DEEP",  'Pie Tin - 8"',                   "1", "178.49", "178.49"),
    @("500CTOUT160",  "Java Box (160oz)",                "2", "94.99",  "189.98"),
    @("245CCGR1914",  "Cake Board - 1/2 Sheet",          "1", "37.99",  "37.99"),
    @("150BB6218N",   "Bag Paper - 6.5x17.75 Window",    "2", "104.99", "209.98"),
    @("150BB4224N",   "Bag Paper - Baguette",            "2", "118.99", "237.98"),
    @("433qlinerbl",  "Sheet Pan Liner - White",         "5", "43.99",  "219.95"),
    @("130TONG6BLK",  'Tong - 6.25" (Black)',            "1", "22.99",  "22.99"),
    @("707U72SPRDBK", 'Black Plastic Spreader 7.51"',    "1", "18.49",  "18.49"),
    @("43302CUPC250", "Container - Muffin (2 Pack)",     "1", "47.99",  "47.99")
)

$startRow = 9
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $rowRange = $ws.Range("A" + $r + ":E" + $r)

    # Temporarily mark the row as Text so numeric-looking values (quantity,
    # cost) are written as literal text rather than being auto-coerced to
    # numbers - matching every other data row in this sheet.
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]

    # The stored value stays text even once the cell format reverts, so
    # restore the default "Normal" style to avoid leaving a stray explicit
    # Text number format on these cells (rows 1-8 have no explicit style).
    $rowRange.Style = "Normal"
}
